$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 (cId 2578889968 / sldId 26415) - "Content Placeholder 2":
# add a new bold red follow-up bullet after the existing "Document Next
# Steps" feedback line.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
$tr3.InsertAfter("`r-> Looking forward for review and feedback from working group.") | Out-Null

# ---------------------------------------------------------------------------
# Slide 4 (cId 269908607 / sldId 2145706288) - "Content Placeholder 2":
# tighten the spacing above the three section headings, reword two
# "Conversation with" call-outs to "Review with", and add the same
# follow-up bullet used on slide 3.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

# "Deployment Next Steps" / "Document Status" / "Document Next Steps"
# headings: spcBef 600 -> 300 (i.e. 6pt -> 3pt).
$tr4.Paragraphs(3, 1).ParagraphFormat.SpaceBefore = 3
$tr4.Paragraphs(5, 1).ParagraphFormat.SpaceBefore = 3
$tr4.Paragraphs(8, 1).ParagraphFormat.SpaceBefore = 3

# "-> Conversation with " -> "-> Review with "
$tr4.Paragraphs(9, 1).Runs(1, 1).Text = "-> Review with "

# "Conversation with " -> "Review with "
$tr4.Paragraphs(10, 1).Runs(1, 1).Text = "Review with "

# Append the same closing bullet that slide 3 received.
$tr4.InsertAfter("`r-> Looking forward for review and feedback from working group.") | Out-Null

# ---------------------------------------------------------------------------
# Slide 6 (cId 559060690 / sldId 2145706289) - "Google Shape;221;p8":
# drop the trailing "Repository: <link>" line and shrink the box to fit
# the now single-line text.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item("Google Shape;221;p8")
$tr6 = $shp6.TextFrame.TextRange
$len6 = $tr6.Length
$tail6 = $tr6.Characters(13, $len6 - 12)
$tail6.Text = " relate to existing service and network topology YANG modules to enable topology visualization."
$shp6.Height = 36.3492125984252
